# "small fix to xlsx" - correct one data row's BLEU/PINC figures and its
# label, and tidy up the leftover UI view-state (selection/scroll) that
# Excel recorded when the file was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data correction on row 7 ---------------------------------------
$ws.Range("D7").Value = 15.69
$ws.Range("E7").Value = 52.94
$ws.Range("F7").Value = "16plays_36LM_mert.2"

# --- chart: the data-label position override is no longer wanted ----
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$series = $chart.SeriesCollection(1)
$dLbls = $series.DataLabels()
$dLbls.Position = -4142

# --- sheet view: select D7, and let the scroll position reset -------
$ws.Range("D7").Select()
